$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-unused chapter rows (14.1, 14.2)
$ws.Rows("22:23").Delete()

# Insert a new row for chapter 19.4 (after 19.1/19.2, before 23.1)
$ws.Rows("27:27").Insert()
$ws.Range("A27").Value = "19.4"
$ws.Range("N27").Formula = "=IF(SUM(B27:M27)>0,1,0)"

# Mark chapters 15.1 and 15.2 (now rows 22/23) as used by lab05 (column D)
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 1

# Update the active selection
$ws.Range("Q23").Select() | Out-Null

# Re-apply the conditional formatting on the shrunk range so dxf/sqref regenerate cleanly
$ws.Cells.FormatConditions.Delete()
$newfc = $ws.Range("N3:N34").FormatConditions.Add(8, 3, "=0")
$newfc.Font.Color = 393372
$newfc.Interior.Color = 13551615
